# Realestate Update resale numbers 2023-06-26 19:44
# Appends a new data row (row 79) to the CityResaleNum sheet with the
# resale-number snapshot captured on 2023-06-26 at 19:39:58.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 79

# Columns A-D hold text (date/time/weekday/week-number) values in this
# sheet even though some look numeric (e.g. "26"), so force the cells to
# text format before assigning them - this keeps them as text instead of
# Excel auto-converting them to a date serial / number.
$ws.Range("A79:D79").NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2023-06-26"
$ws.Cells.Item($row, 2).Value = "19:39:58"
$ws.Cells.Item($row, 3).Value = "Monday"
$ws.Cells.Item($row, 4).Value = "26"

# Columns E-T are plain numeric resale-count values.
$ws.Cells.Item($row, 5).Value = 122828
$ws.Cells.Item($row, 6).Value = 134281
$ws.Cells.Item($row, 7).Value = 163681
$ws.Cells.Item($row, 8).Value = 133613
$ws.Cells.Item($row, 9).Value = 177264
$ws.Cells.Item($row, 10).Value = 115096
$ws.Cells.Item($row, 11).Value = 203598
$ws.Cells.Item($row, 12).Value = 226259
$ws.Cells.Item($row, 13).Value = 176279
$ws.Cells.Item($row, 14).Value = 104293
$ws.Cells.Item($row, 15).Value = 39657
$ws.Cells.Item($row, 16).Value = 33783
$ws.Cells.Item($row, 17).Value = 52197
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 35994
$ws.Cells.Item($row, 20).Value = -1
